# Applies the "PEPified the code, trying to follow established conventions
# more" commit: clears leftover test/placeholder strings ("jalla", "Hei",
# "AEring", "the ", "Laerling", "ehi", "Please input values") out of the
# data tables - replacing them with real numeric data - appends a new data
# row to the "Material Removal Rate" sheet, and moves the active
# tab/selection from "Surface Roughness" back to "Cutting Speed".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Cutting Speed"
# ---------------------------------------------------------------------
$cuttingSpeed = $wb.Worksheets.Item("Cutting Speed")

$cuttingSpeed.Range("B5").Value = 0

$cuttingSpeed.Range("B6").Value = 0
$cuttingSpeed.Range("C6").Value = 80
$cuttingSpeed.Range("D6").Value = 6
$cuttingSpeed.Range("E6").Value = 0.12

$cuttingSpeed.Range("C7").Value = 0
$cuttingSpeed.Range("D7").Value = 0
$cuttingSpeed.Range("E7").Value = 0

$cuttingSpeed.Range("B10").Value = 210
$cuttingSpeed.Range("C10").Value = 80
$cuttingSpeed.Range("D10").Value = 6
$cuttingSpeed.Range("E10").Value = 0.12

# ---------------------------------------------------------------------
# Sheet "Material Removal Rate"
# ---------------------------------------------------------------------
$mrr = $wb.Worksheets.Item("Material Removal Rate")

$mrr.Range("B16").Value = 0
$mrr.Range("E16").Value = 0

$mrr.Range("E18").Value = 0

# New row of real data appended at the bottom of the table
$mrr.Range("B26").Value = 1
$mrr.Range("C26").Value = 60
$mrr.Range("D26").Value = 602
$mrr.Range("E26").Value = 36.12
$mrr.Range("F26").Value = "cm³/min"

# ---------------------------------------------------------------------
# Sheet "Helix Angle"
# ---------------------------------------------------------------------
$helixAngle = $wb.Worksheets.Item("Helix Angle")

$helixAngle.Range("E6").Value = 0

# ---------------------------------------------------------------------
# Sheet "Ramp Angle"
# ---------------------------------------------------------------------
$rampAngle = $wb.Worksheets.Item("Ramp Angle")

$rampAngle.Range("D7").Value = 0

# ---------------------------------------------------------------------
# Selections per sheet (also drives which sheet ends up "tabSelected"
# and the workbook's activeTab) - activating each sheet and selecting a
# cell on it records that sheet's last selection; activating the sheet
# we want as the final active tab LAST makes it the one that is
# persisted as tabSelected/activeTab.
# ---------------------------------------------------------------------
$mrr.Activate()
$mrr.Range("C29").Select()

$helixAngle.Activate()
$helixAngle.Range("E7").Select()

$rampAngle.Activate()
$rampAngle.Range("D8").Select()

$surfaceRoughness = $wb.Worksheets.Item("Surface Roughness")
$surfaceRoughness.Activate()
$surfaceRoughness.Range("C12").Select()

# "Cutting Speed" is activated last so it becomes the active tab, with
# the selection sitting on D11.
$cuttingSpeed.Activate()
$cuttingSpeed.Range("D11").Select()
